$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift the existing data rows 2-7 down to 3-8 (direct literal writes, no
# real row insert, so the formatting placeholder rows further down stay put).

# Row 8 (was row 7): 18-2-8 / AH317 Royal
$ws.Range("A8").Value = "18-2-8"
$ws.Range("B8").Value = "AH317 Royal"
$ws.Range("C8").Value = 23
$ws.Range("D8").Value = 150
$ws.Range("E8").Value = "30/1/23"
$ws.Range("F8").Value = "SZ"

# Row 7 (was row 6): 18-2-8 / AH230 Red
$ws.Range("A7").Value = "18-2-8"
$ws.Range("B7").Value = "AH230 Red"
$ws.Range("C7").Value = 12
$ws.Range("D7").Value = 150
$ws.Range("E7").Value = "30/1/23"
$ws.Range("F7").Value = "SZ"

# Row 6 (was row 5): 18-2-4 / Ah256 Red
$ws.Range("A6").Value = "18-2-4"
$ws.Range("B6").Value = "Ah256 Red"
$ws.Range("C6").Value = 12
$ws.Range("D6").Value = 150
$ws.Range("E6").Value = "30/1/23"
$ws.Range("F6").Value = "SZ"

# Row 5 (was row 4): 18-2-4 / AH256 Black
$ws.Range("A5").Value = "18-2-4"
$ws.Range("B5").Value = "AH256 Black"
$ws.Range("C5").Value = 5
$ws.Range("D5").Value = 150
$ws.Range("E5").Value = "30/1/23"
$ws.Range("F5").Value = "SZ"

# Row 4 (was row 3): 18-2-3 / G4152 Grey/Black, Pieces now 15
$ws.Range("A4").Value = "18-2-3"
$ws.Range("B4").Value = "G4152 Grey/Black"
$ws.Range("C4").Value = 12
$ws.Range("D4").Value = 15
$ws.Range("E4").Value = "30/1/23"
$ws.Range("F4").Value = "SZ"
# B column carries the odd date-style formatting that travelled with this
# item's name cell - copy it across instead of the (now stale) source row.
$ws.Range("B3").Copy()
$ws.Range("B4").PasteSpecial(-4122)

# Row 3 (was row 2): now only the location code remains
$ws.Range("A3").Value = "12-2-9"
$ws.Range("B3:F3").Clear()

# Row 2 (new): a new pallet/location entry
$ws.Range("A2").Value = "1A-3-1"

# Restore the active selection to match the saved view state
$ws.Range("N23").Select()
